$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Chequeo de directorios de HL y HE60: actualizar el rango de Id_run.
$ws.Range("C2").Value = 130
$ws.Range("D2").Value = 150
